$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("H2").Value = 43
$ws.Range("I3").Value = 75
$ws.Range("J3").Value = 69
$ws.Range("E9").Value = 173
$ws.Range("F9").Value = 198
$ws.Range("I9").Value = 211
$ws.Range("B10").Value = 408
$ws.Range("C10").Value = 493
$ws.Range("D10").Value = 664
$ws.Range("E10").Value = 739
$ws.Range("F10").Value = 859
$ws.Range("G10").Value = 505
$ws.Range("H10").Value = 168
$ws.Range("I10").Value = 309
$ws.Range("J10").Value = 271
$ws.Range("B11").Value = 601
$ws.Range("C11").Value = 734
$ws.Range("D11").Value = 921
$ws.Range("E11").Value = 998
$ws.Range("F11").Value = 1139
$ws.Range("G11").Value = 791
$ws.Range("H11").Value = 402
$ws.Range("I11").Value = 638
$ws.Range("J11").Value = 562

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("I5").Value = 5
$ws.Range("I7").Value = 10

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("E6").Value = 11
$ws.Range("E8").Value = 44

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 14
$ws.Range("F7").Value = 25
$ws.Range("J7").Value = 26
$ws.Range("D8").Value = 197
$ws.Range("E8").Value = 218
$ws.Range("F8").Value = 261
$ws.Range("H8").Value = 18
$ws.Range("D9").Value = 232
$ws.Range("E9").Value = 257
$ws.Range("F9").Value = 292
$ws.Range("H9").Value = 48
$ws.Range("J9").Value = 98

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("F6").Value = 47
$ws.Range("F7").Value = 53

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("F6").Value = 17
$ws.Range("F7").Value = 28

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("C8").Value = 37
$ws.Range("D10").Value = 17
$ws.Range("I21").Value = 10
$ws.Range("E22").Value = 7
$ws.Range("B28").Value = 38
$ws.Range("H28").Value = 27
$ws.Range("E32").Value = 44
$ws.Range("B38").Value = 2
$ws.Range("G43").Value = 7
$ws.Range("G49").Value = 2
$ws.Range("I49").Value = 5
$ws.Range("F50").Value = 28
$ws.Range("J52").Value = 16
$ws.Range("D53").Value = 232
$ws.Range("E53").Value = 257
$ws.Range("F53").Value = 292
$ws.Range("H53").Value = 48
$ws.Range("J53").Value = 98
$ws.Range("C61").Value = 14
$ws.Range("D62").Value = 13
$ws.Range("C67").Value = 4
$ws.Range("F70").Value = 53
$ws.Range("F72").Value = 21
$ws.Range("I76").Value = 14
$ws.Range("C77").Value = 29
$ws.Range("B82").Value = 11
$ws.Range("B83").Value = 3
$ws.Range("B98").Value = 601
$ws.Range("C98").Value = 734
$ws.Range("D98").Value = 921
$ws.Range("E98").Value = 998
$ws.Range("F98").Value = 1139
$ws.Range("G98").Value = 791
$ws.Range("H98").Value = 402
$ws.Range("I98").Value = 638
$ws.Range("J98").Value = 562

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("H2").Value = 5
$ws.Range("B5").Value = 23
$ws.Range("B6").Value = 38
$ws.Range("H6").Value = 27

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("I7").Value = 9
$ws.Range("I8").Value = 14

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("B5").Value = 6
$ws.Range("B6").Value = 11

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("C5").Value = 4
$ws.Range("C6").Value = 4

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 16

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("D5").Value = 12
$ws.Range("D6").Value = 13

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("C8").Value = 20
$ws.Range("C9").Value = 29

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("I3").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 2
$ws.Range("I6").Value = 5

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("B6").Value = 1
$ws.Range("B7").Value = 3

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("E5").Value = 6
$ws.Range("E6").Value = 7

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("F5").Value = 16
$ws.Range("F6").Value = 21

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("D5").Value = 16
$ws.Range("D6").Value = 17

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 2

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("C7").Value = 17
$ws.Range("C8").Value = 37

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 7
